$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 33166.332
$ws.Range("I18").Value = 1500.0
$ws.Range("K18").Value = 1500.0
$ws.Range("M18").Value = -1216.0
$ws.Range("H32").Value = 6000.5
$ws.Range("J32").Value = 5000.0
$ws.Range("L32").Value = 5000.0
$ws.Range("N32").Value = -5652.0
$ws.Range("H69").Value = 0.0
$ws.Range("I69").Value = 0.0
$ws.Range("K69").Value = 0.0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0.0
$ws.Range("I72").Value = 0.0
$ws.Range("K72").Value = 0.0
$ws.Range("M72").ClearContents()
$ws.Range("H80").Value = 871.625
$ws.Range("I80").Value = 763.55554
$ws.Range("J80").Value = 936.4667
$ws.Range("K80").Value = 2290.66662
$ws.Range("L80").Value = 2809.4001
$ws.Range("M80").Value = -1292.66662
$ws.Range("N80").Value = -4805.4001
$ws.Range("H83").Value = 871.625
$ws.Range("I83").Value = 763.55554
$ws.Range("J83").Value = 936.4667
$ws.Range("K83").Value = 6871.99986
$ws.Range("L83").Value = 8428.2003
$ws.Range("M83").Value = -1879.99986
$ws.Range("N83").Value = -18412.2003
$ws.Range("H106").Value = 9575.929
$ws.Range("I106").Value = 1949.0
$ws.Range("K106").Value = 1949.0
$ws.Range("M106").Value = -1318.0
$ws.Range("H107").Value = 631.3333
$ws.Range("I107").Value = 631.3333
$ws.Range("K107").Value = 631.3333
$ws.Range("M107").Value = 1288.6667
$ws.Range("H138").Value = 2152.6428
$ws.Range("I138").Value = 1686.3636
$ws.Range("K138").Value = 5059.0908
$ws.Range("M138").Value = 80.90920000000006

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 468.5
$ws.Range("I22").Value = 468.5
$ws.Range("K22").Value = 468.5
$ws.Range("M22").Value = -295.5
$ws.Range("H36").Value = 12681.714
$ws.Range("I36").Value = 8900.0
$ws.Range("J36").Value = 15518.0
$ws.Range("K36").Value = 8900.0
$ws.Range("L36").Value = 15518.0
$ws.Range("M36").Value = -8366.0
$ws.Range("N36").Value = -16586.0
$ws.Range("H74").Value = 80000.0
$ws.Range("J74").Value = 80000.0
$ws.Range("L74").Value = 80000.0
$ws.Range("N74").Value = -81872.0
$ws.Range("H77").Value = 80000.0
$ws.Range("J77").Value = 80000.0
$ws.Range("L77").Value = 240000.0
$ws.Range("N77").Value = -249360.0
$ws.Range("H81").Value = 22221.5
$ws.Range("J81").Value = 22221.5
$ws.Range("L81").Value = 22221.5
$ws.Range("N81").Value = -24343.5
$ws.Range("H84").Value = 22221.5
$ws.Range("J84").Value = 22221.5
$ws.Range("L84").Value = 66664.5
$ws.Range("N84").Value = -77272.5
$ws.Range("H134").Value = 2712.8096
$ws.Range("I134").Value = 2429.9285
$ws.Range("K134").Value = 7289.7855
$ws.Range("M134").Value = -4754.7855
$ws.Range("H139").Value = 80000.0
$ws.Range("J139").Value = 80000.0
$ws.Range("L139").Value = 80000.0
$ws.Range("N139").Value = -90280.0
$ws.Range("H140").Value = 94113.0
$ws.Range("J140").Value = 94113.0
$ws.Range("L140").Value = 94113.0
$ws.Range("N140").Value = -104473.0
$ws.Range("H141").Value = 95990.6
$ws.Range("J141").Value = 95990.6
$ws.Range("L141").Value = 95990.6
$ws.Range("N141").Value = -106350.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 1422.0
$ws.Range("J17").Value = 1422.0
$ws.Range("L17").Value = 1422.0
$ws.Range("N17").Value = -1770.0
$ws.Range("H19").Value = 1487.6666
$ws.Range("I19").Value = 419.375
$ws.Range("J19").Value = 3624.25
$ws.Range("K19").Value = 419.375
$ws.Range("L19").Value = 3624.25
$ws.Range("M19").Value = -249.375
$ws.Range("N19").Value = -3964.25
$ws.Range("H23").Value = 9.0
$ws.Range("I23").Value = 9.0
$ws.Range("J23").Value = 0.0
$ws.Range("K23").Value = 9.0
$ws.Range("L23").Value = 0.0
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = 231.0
$ws.Range("H24").Value = 1487.6666
$ws.Range("I24").Value = 419.375
$ws.Range("J24").Value = 3624.25
$ws.Range("K24").Value = 419.375
$ws.Range("L24").Value = 3624.25
$ws.Range("M24").Value = -249.375
$ws.Range("N24").Value = -3964.25
$ws.Range("H25").Value = 805.5
$ws.Range("J25").Value = 0.0
$ws.Range("L25").Value = 0.0
$ws.Range("N25").ClearContents()
$ws.Range("H27").Value = 9.0
$ws.Range("I27").Value = 9.0
$ws.Range("J27").Value = 0.0
$ws.Range("K27").Value = 9.0
$ws.Range("L27").Value = 0.0
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = 183.0

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1670.5555
$ws.Range("I3").Value = 1670.5555
$ws.Range("K3").Value = 5011.666499999999
$ws.Range("M3").Value = -4899.666499999999
$ws.Range("H38").Value = 225.8
$ws.Range("J38").Value = 450.0
$ws.Range("L38").Value = 1350.0
$ws.Range("N38").Value = -2044.0
$ws.Range("H122").Value = 829.45
$ws.Range("J122").Value = 899.94446
$ws.Range("L122").Value = 8099.50014
$ws.Range("N122").Value = -12999.50014
$ws.Range("H132").Value = 1550.0
$ws.Range("I132").Value = 1100.0
$ws.Range("K132").Value = 9900.0
$ws.Range("M132").Value = -7370.0

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1559.2354
$ws.Range("I22").Value = 1547.6
$ws.Range("K22").Value = 1547.6
$ws.Range("M22").Value = -1252.6
$ws.Range("H27").Value = 1559.2354
$ws.Range("I27").Value = 1547.6
$ws.Range("K27").Value = 1547.6
$ws.Range("M27").Value = -1440.6
$ws.Range("H61").Value = 13810.8
$ws.Range("I61").Value = 18030.572
$ws.Range("J61").Value = 3964.6667
$ws.Range("K61").Value = 18030.572
$ws.Range("L61").Value = 3964.6667
$ws.Range("M61").Value = -17828.572
$ws.Range("N61").Value = -4368.6667
$ws.Range("H82").Value = 1128.1111
$ws.Range("J82").Value = 1090.2
$ws.Range("L82").Value = 1090.2
$ws.Range("N82").Value = -1812.2
$ws.Range("H85").Value = 1128.1111
$ws.Range("J85").Value = 1090.2
$ws.Range("L85").Value = 1090.2
$ws.Range("N85").Value = -3586.2
$ws.Range("H113").Value = 13810.8
$ws.Range("I113").Value = 18030.572
$ws.Range("J113").Value = 3964.6667
$ws.Range("K113").Value = 18030.572
$ws.Range("L113").Value = 3964.6667
$ws.Range("M113").Value = -15860.572
$ws.Range("N113").Value = -8304.6667
